$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force Text format on the data range so that numeric-looking
# strings (e.g. "560.80", "1.00") are not auto-converted to numbers by Excel.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '64.115.07'
$ws.Range("E2").Value = '  +1.04%  '

$ws.Range("D3").Value = '3.091.95'
$ws.Range("E3").Value = '  +0.79%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").Value = '560.80'
$ws.Range("E5").Value = '  +2.10%  '

$ws.Range("D6").Value = '144.69'
$ws.Range("E6").Value = '  +2.90%  '

$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.07%  '

$ws.Range("D8").Value = '3.088.88'
$ws.Range("E8").Value = '  +0.88%  '

$ws.Range("D9").Value = '0.505'
$ws.Range("E9").Value = '  +0.66%  '

$ws.Range("D10").Value = '0.153'
$ws.Range("E10").Value = '  +1.70%  '

$ws.Range("D11").Value = '6.15'
$ws.Range("E11").Value = '  -5.95%  '

$ws.Range("D12").Value = '0.471'
$ws.Range("E12").Value = '  +3.62%  '

$ws.Range("E13").Value = '  +0.24%  '

$ws.Range("D14").Value = '35.13'
$ws.Range("E14").Value = '  +0.78%  '

$ws.Range("D15").Value = '3.594.83'
$ws.Range("E15").Value = '  +0.91%  '

$ws.Range("D16").Value = '64.171.56'
$ws.Range("E16").Value = '  +1.23%  '

$ws.Range("D17").Value = '3.090.03'
$ws.Range("E17").Value = '  +0.78%  '

$ws.Range("E18").Value = '  +1.39%  '

$ws.Range("D19").Value = '6.76'
$ws.Range("E19").Value = '  +0.19%  '

$ws.Range("D20").Value = '482.83'
$ws.Range("E20").Value = '  +0.20%  '

$ws.Range("D21").Value = '13.96'
$ws.Range("E21").Value = '  +2.04%  '

$ws.Range("D22").Value = '0.674'
$ws.Range("E22").Value = '  +0.28%  '

$ws.Range("E23").Value = '  +3.89%  '

$ws.Range("D24").Value = '14.03'
$ws.Range("E24").Value = '  +11.05%  '

$ws.Range("D25").Value = '81.21'
$ws.Range("E25").Value = '  +0.55%  '

$ws.Range("E26").Value = '  +0.04%  '

$ws.Range("D27").Value = '2.80'
$ws.Range("E27").Value = '  +1.68%  '

$ws.Range("D28").Value = '8.00'
$ws.Range("E28").Value = '  +0.95%  '

$ws.Range("E29").Value = '  +3.29%  '

$ws.Range("E30").Value = '  +0.01%  '

$ws.Range("D31").Value = '26.32'
$ws.Range("E31").Value = '  +0.73%  '

$ws.Range("D32").Value = '1.15'
$ws.Range("E32").Value = '  +0.18%  '

$ws.Range("E33").Value = '  +0.75%  '

$ws.Range("D34").Value = '5.58'
$ws.Range("E34").Value = '  -2.56%  '

$ws.Range("B35").Value = 'OKB'
$ws.Range("C35").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D35").Value = '55.82'
$ws.Range("E35").Value = '  +0.43%  '

$ws.Range("B36").Value = 'Filecoin'
$ws.Range("C36").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D36").Value = '6.19'
$ws.Range("E36").Value = '  +3.46%  '

$ws.Range("D37").Value = '452.96'
$ws.Range("E37").Value = '  -3.22%  '

$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '0.0408'
$ws.Range("E38").Value = '  +2.87%  '

$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").Value = '0.0820'
$ws.Range("E39").Value = '  -0.11%  '

$ws.Range("B40").Value = 'dogwifhat'
$ws.Range("C40").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D40").Value = '2.91'
$ws.Range("E40").Value = '  +13.02%  '

$ws.Range("D41").Value = '2.983.08'
$ws.Range("E41").Value = '  -2.87%  '

$ws.Range("D42").Value = '8.25'
$ws.Range("E42").Value = '  -0.12%  '

$ws.Range("E43").Value = '  -3.69%  '

$ws.Range("D44").Value = '27.91'
$ws.Range("E44").Value = '  -0.26%  '

$ws.Range("D45").Value = '0.260'
$ws.Range("E45").Value = '  +2.32%  '

$ws.Range("B46").Value = 'USDe'
$ws.Range("C46").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D46").Value = '1.00'
$ws.Range("E46").Value = '  -0.01%  '

$ws.Range("B47").Value = 'Fetch.AI'
$ws.Range("C47").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D47").Value = '2.13'
$ws.Range("E47").Value = '  +3.41%  '

$ws.Range("E48").Value = '  +1.93%  '

$ws.Range("D49").Value = '120.11'
$ws.Range("E49").Value = '  +3.14%  '

$ws.Range("D50").Value = '0.0₃0514'
$ws.Range("E50").Value = '  +0.83%  '

$ws.Range("E51").Value = '  +0.19%  '

# Restore original (default) cell formatting now that the text values are set.
$dataRange.ClearFormats()